$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checks")

# Remove the autofilter criteria that was limiting the "Trait" column (C)
# to "COSTR" only. This unhides every previously-filtered-out row
# (rows 20-58) while keeping the AutoFilter/sort state in place - i.e.
# calculating the number of checks per genotype across ALL traits, not
# just COSTR.
$ws.ShowAllData()

# Underline the trait value in C2 (now visually distinguishing the
# first/representative trait row after the filter was cleared).
$ws.Range("C2").Font.Underline = 2
